# Scheduled-runner price/profit refresh across Leve sheets.
# Updates currentAveragePrice(NQ/HQ) / LevePrice(NQ/HQ) / LeveProfit(NQ/HQ)
# columns (H:N) for the rows whose market data changed since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 950
$ws.Range("I18").Value = 950
$ws.Range("K18").Value = 950
$ws.Range("M18").Value = -666

$ws.Range("H51").Value = 9623.666999999999
$ws.Range("I51").Value = 9998.5
$ws.Range("J51").Value = 9548.700000000001
$ws.Range("K51").Value = 9998.5
$ws.Range("L51").Value = 9548.700000000001
$ws.Range("M51").Value = -9514.5
$ws.Range("N51").Value = -10516.7

$ws.Range("H80").Value = 1829
$ws.Range("I80").Value = 1667.3334
$ws.Range("J80").Value = 1950.25
$ws.Range("K80").Value = 5002.0002
$ws.Range("L80").Value = 5850.75
$ws.Range("M80").Value = -4004.0002
$ws.Range("N80").Value = -7846.75

$ws.Range("H83").Value = 1829
$ws.Range("I83").Value = 1667.3334
$ws.Range("J83").Value = 1950.25
$ws.Range("K83").Value = 15006.0006
$ws.Range("L83").Value = 17552.25
$ws.Range("M83").Value = -10014.0006
$ws.Range("N83").Value = -27536.25

$ws.Range("H112").Value = 3833.3333
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 5500
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 16500
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -18716

$ws.Range("H132").Value = 6383.857
$ws.Range("I132").Value = 6781.1665
$ws.Range("K132").Value = 20343.4995
$ws.Range("M132").Value = -17813.4995

$ws.Range("H138").Value = 1957.1428
$ws.Range("I138").Value = 1300
$ws.Range("J138").Value = 2450
$ws.Range("K138").Value = 3900
$ws.Range("L138").Value = 7350
$ws.Range("M138").Value = 1240
$ws.Range("N138").Value = -17630

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16125
$ws.Range("I32").Value = 16125
$ws.Range("K32").Value = 16125
$ws.Range("M32").Value = -15838

$ws.Range("H110").Value = 723.6
$ws.Range("I110").Value = 723.6
$ws.Range("K110").Value = 723.6
$ws.Range("M110").Value = 1321.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3437
$ws.Range("I86").Value = 3437
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3437
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2314
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 3437
$ws.Range("I89").Value = 3437
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 17185
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -11569
$ws.Range("N89").ClearContents()

$ws.Range("H99").Value = 3150
$ws.Range("I99").Value = 3150
$ws.Range("K99").Value = 3150
$ws.Range("M99").Value = -1652

$ws.Range("H105").Value = 20992.5
$ws.Range("I105").Value = 31807
$ws.Range("J105").Value = 2968.3333
$ws.Range("K105").Value = 31807
$ws.Range("L105").Value = 2968.3333
$ws.Range("M105").Value = -30060
$ws.Range("N105").Value = -6462.3333

$ws.Range("H123").Value = 100000
$ws.Range("J123").Value = 100000
$ws.Range("L123").Value = 100000
$ws.Range("N123").Value = -109800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1154.4
$ws.Range("I31").Value = 1154.4
$ws.Range("K31").Value = 1154.4
$ws.Range("M31").Value = -859.4000000000001

$ws.Range("H34").Value = 1154.4
$ws.Range("I34").Value = 1154.4
$ws.Range("K34").Value = 1154.4
$ws.Range("M34").Value = -952.4000000000001

$ws.Range("H69").Value = 3125
$ws.Range("I69").Value = 3125
$ws.Range("K69").Value = 3125
$ws.Range("M69").Value = -2376

$ws.Range("H72").Value = 3125
$ws.Range("I72").Value = 3125
$ws.Range("K72").Value = 9375
$ws.Range("M72").Value = -5631

$ws.Range("H76").Value = 100
$ws.Range("I76").Value = 100
$ws.Range("K76").Value = 100
$ws.Range("M76").Value = 215

$ws.Range("H79").Value = 100
$ws.Range("I79").Value = 100
$ws.Range("K79").Value = 100
$ws.Range("M79").Value = 992

$ws.Range("H88").Value = 29885.75
$ws.Range("J88").Value = 29885.75
$ws.Range("L88").Value = 29885.75
$ws.Range("N88").Value = -30697.75

$ws.Range("H91").Value = 29885.75
$ws.Range("J91").Value = 29885.75
$ws.Range("L91").Value = 29885.75
$ws.Range("N91").Value = -32693.75

$ws.Range("H95").Value = 6363.8
$ws.Range("J95").Value = 6363.8
$ws.Range("L95").Value = 6363.8
$ws.Range("N95").Value = -11855.8

$ws.Range("H96").Value = 12131.714
$ws.Range("J96").Value = 12131.714
$ws.Range("L96").Value = 12131.714
$ws.Range("N96").Value = -17623.714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 900
$ws.Range("J75").Value = 900
$ws.Range("L75").Value = 2700
$ws.Range("N75").Value = -4696

$ws.Range("H78").Value = 900
$ws.Range("J78").Value = 900
$ws.Range("L78").Value = 8100
$ws.Range("N78").Value = -18084

$ws.Range("H80").Value = 8250
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 15000
$ws.Range("K80").Value = 4500
$ws.Range("L80").Value = 45000
$ws.Range("M80").Value = -3564
$ws.Range("N80").Value = -46872

$ws.Range("H83").Value = 8250
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 15000
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 135000
$ws.Range("M83").Value = -8820
$ws.Range("N83").Value = -144360

$ws.Range("H134").Value = 1250
$ws.Range("I134").Value = 1250
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3750
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 1320
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 866
$ws.Range("I46").Value = 799.5
$ws.Range("J46").Value = 999
$ws.Range("K46").Value = 799.5
$ws.Range("L46").Value = 999
$ws.Range("M46").Value = -611.5
$ws.Range("N46").Value = -1375

$ws.Range("H61").Value = 3986
$ws.Range("I61").Value = 3817.6667
$ws.Range("J61").Value = 4238.5
$ws.Range("K61").Value = 3817.6667
$ws.Range("L61").Value = 4238.5
$ws.Range("M61").Value = -3615.6667
$ws.Range("N61").Value = -4642.5

$ws.Range("H113").Value = 3986
$ws.Range("I113").Value = 3817.6667
$ws.Range("J113").Value = 4238.5
$ws.Range("K113").Value = 3817.6667
$ws.Range("L113").Value = 4238.5
$ws.Range("M113").Value = -1647.6667
$ws.Range("N113").Value = -8578.5

$ws.Range("H122").Value = 4317.909
$ws.Range("I122").Value = 4166.4443
$ws.Range("K122").Value = 12499.3329
$ws.Range("M122").Value = -10049.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 15590
$ws.Range("J69").Value = 15590
$ws.Range("L69").Value = 15590
$ws.Range("N69").Value = -17088

$ws.Range("H72").Value = 15590
$ws.Range("J72").Value = 15590
$ws.Range("L72").Value = 46770
$ws.Range("N72").Value = -54258
